$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline approach: for values that Excel would otherwise
# auto-convert to a number (plain numeric-looking text), force the cell
# to Text format first, assign the string, then reset the style back to
# "Normal" so no stray number-format/style is left behind.

$ws.Range('D2').Value = '79.727.10'
$ws.Range('E2').Value = '  +4.45%  '
$ws.Range('D3').Value = '3.200.19'
$ws.Range('E3').Value = '  +5.29%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '641.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.246'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +19.06%  '
$ws.Range('E9').Value = '  +11.43%  '
$ws.Range('D10').Value = '3.197.04'
$ws.Range('E10').Value = '  +5.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.611'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +39.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000265'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +37.32%  '
$ws.Range('E13').Value = '  +3.59%  '
$ws.Range('E14').Value = '  +3.27%  '
$ws.Range('D15').Value = '3.789.90'
$ws.Range('E15').Value = '  +5.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +13.15%  '
$ws.Range('D17').Value = '79.622.50'
$ws.Range('E17').Value = '  +4.32%  '
$ws.Range('D18').Value = '3.208.07'
$ws.Range('E18').Value = '  +5.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '443.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +17.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.98'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +28.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +22.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.85'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.98%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.367.57'
$ws.Range('E25').Value = '  +5.28%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '77.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.86%  '
$ws.Range('E27').Value = '  +13.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  +16.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.20'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +11.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.56'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +11.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '552.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.156'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +33.20%  '
$ws.Range('E35').Value = '  +6.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +13.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.123'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +18.30%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.420'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.22%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '165.50'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '194.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  +11.49%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.83'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +11.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +12.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.808'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.96%  '
$ws.Range('E48').Value = '  +6.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.98'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.24'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +18.73%  '
$ws.Range('E51').Value = '  +7.77%  '
